$wb = $excel.ActiveWorkbook

# --- 1. Update shared-string text: "Ready for handoff" -> "In Translation" ---
# This string appears in the Status column in all three sheets.
$ws1 = $wb.Worksheets.Item(1)   # "Overview"
$ws2 = $wb.Worksheets.Item(2)   # "zh-cn"
$ws3 = $wb.Worksheets.Item(3)   # "de-de"

$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws2.Range("C2").Value = "In Translation"
$ws3.Range("C2").Value = "In Translation"

# --- 2. Shrink the "Status" columns to their new (narrower) width ---
# Original stored width 17.2159881591797 characters -> new stored width
# 13.4101845877511 characters. Excel quantizes ColumnWidth to whole pixels
# (MDW grid), so we use the ColumnWidth value that lands closest to the
# target stored width.
$newColumnWidth = 12.5

$ws1.Range("E1").ColumnWidth = $newColumnWidth
$ws1.Range("F1").ColumnWidth = $newColumnWidth
$ws2.Range("C1").ColumnWidth = $newColumnWidth
$ws3.Range("C1").ColumnWidth = $newColumnWidth
